$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 497, pushing the existing rows 497:532 down to 498:533.
$ws.Rows.Item(497).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(497, 1).Value = 5
$ws.Cells.Item(497, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(497, 3).Value = "Maule"
$ws.Cells.Item(497, 4).Value = 45021
$ws.Cells.Item(497, 5).Value = 7
$ws.Cells.Item(497, 6).Value = 100112032
$ws.Cells.Item(497, 7).Value = "Zapallo italiano"
$ws.Cells.Item(497, 8).Value = "Sin especificar"
$ws.Cells.Item(497, 9).Value = "Primera"
$ws.Cells.Item(497, 10).Value = 300
$ws.Cells.Item(497, 11).Value = 5000
$ws.Cells.Item(497, 12).Value = 5000
$ws.Cells.Item(497, 13).Value = 5000
$ws.Cells.Item(497, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(497, 15).Value = "Región del Maule"
$ws.Cells.Item(497, 16).Value = 100
$ws.Cells.Item(497, 17).Value = 50
$ws.Cells.Item(497, 18).Value = "Hortaliza"

# Column D carries a date number format on every data row; make sure the
# newly inserted row's date cell keeps the same formatting.
$ws.Cells.Item(497, 4).NumberFormat = $ws.Cells.Item(498, 4).NumberFormat
